$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 8.658781320990508
$ws.Cells.Item(2, 3).Value = 5.776349252250968
$ws.Cells.Item(2, 4).Value = 5.210700280086564
$ws.Cells.Item(2, 5).Value = 12.80913531481467
$ws.Cells.Item(2, 6).Value = 25.85055980118697
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 11).Value = 7.881343650347793
$ws.Cells.Item(2, 13).Value = 13.37592750253274
$ws.Cells.Item(2, 15).Value = 23.12837045790585
$ws.Cells.Item(3, 2).Value = 8.375255007433289
$ws.Cells.Item(3, 3).Value = 5.677267269457206
$ws.Cells.Item(3, 4).Value = 5.168341004403166
$ws.Cells.Item(3, 5).Value = 12.59790147270122
$ws.Cells.Item(3, 6).Value = 25.85822056369299
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 11).Value = 7.690156342032529
$ws.Cells.Item(3, 13).Value = 13.20665415481008
$ws.Cells.Item(3, 15).Value = 23.18463731653483
$ws.Cells.Item(4, 2).Value = 8.197753954047911
$ws.Cells.Item(4, 3).Value = 5.614898613607278
$ws.Cells.Item(4, 4).Value = 5.141725974095148
$ws.Cells.Item(4, 5).Value = 12.47089036794712
$ws.Cells.Item(4, 6).Value = 25.86990825973225
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 11).Value = 7.57160353157494
$ws.Cells.Item(4, 13).Value = 13.10492469694338
$ws.Cells.Item(4, 15).Value = 23.2243900530789
$ws.Cells.Item(5, 2).Value = 8.124679888736626
$ws.Cells.Item(5, 3).Value = 5.589116164458923
$ws.Cells.Item(5, 4).Value = 5.130731789762076
$ws.Cells.Item(5, 5).Value = 12.41988153139805
$ws.Cells.Item(5, 6).Value = 25.87642538715114
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 11).Value = 7.523071554954678
$ws.Cells.Item(5, 13).Value = 13.06407271075921
$ws.Cells.Item(5, 15).Value = 23.24189489959304
$ws.Cells.Item(6, 2).Value = 8.112504868418458
$ws.Cells.Item(6, 3).Value = 5.584813414568559
$ws.Cells.Item(6, 4).Value = 5.128897375259002
$ws.Cells.Item(6, 5).Value = 12.41145892354471
$ws.Cells.Item(6, 6).Value = 25.87761343013088
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 11).Value = 7.515001680833175
$ws.Cells.Item(6, 13).Value = 13.05732713711795
$ws.Cells.Item(6, 15).Value = 23.24488029575637
$ws.Cells.Item(7, 2).Value = 8.196771289259445
$ws.Cells.Item(7, 3).Value = 5.614552362320956
$ws.Cells.Item(7, 4).Value = 5.141578297866685
$ws.Cells.Item(7, 5).Value = 12.47019931532172
$ws.Cells.Item(7, 6).Value = 25.86998905223752
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 11).Value = 7.570949809122583
$ws.Cells.Item(7, 13).Value = 13.10437124483398
$ws.Cells.Item(7, 15).Value = 23.22462084931497
$ws.Cells.Item(8, 2).Value = 8.561797887944978
$ws.Cells.Item(8, 3).Value = 5.742514597091821
$ws.Cells.Item(8, 4).Value = 5.196223074784291
$ws.Cells.Item(8, 5).Value = 12.73578857916102
$ws.Cells.Item(8, 6).Value = 25.85175090478648
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 11).Value = 7.815705683912501
$ws.Cells.Item(8, 13).Value = 13.31713419112213
$ws.Cells.Item(8, 15).Value = 23.14668902793693
$ws.Cells.Item(9, 2).Value = 9.245829622695178
$ws.Cells.Item(9, 3).Value = 5.980520859281286
$ws.Cells.Item(9, 4).Value = 5.298395032029184
$ws.Cells.Item(9, 5).Value = 13.27465665702027
$ws.Cells.Item(9, 6).Value = 25.8714380582451
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 11).Value = 8.36248194742536
$ws.Cells.Item(9, 13).Value = 13.74965334870995
$ws.Cells.Item(9, 15).Value = 23.03530757062607
$ws.Cells.Item(10, 2).Value = 9.723405655030964
$ws.Cells.Item(10, 3).Value = 6.146528865376228
$ws.Cells.Item(10, 4).Value = 5.370192279400364
$ws.Cells.Item(10, 5).Value = 13.67717495384454
$ws.Cells.Item(10, 6).Value = 25.91969059850067
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 11).Value = 8.810338483946481
$ws.Cells.Item(10, 13).Value = 14.07375893708159
$ws.Cells.Item(10, 15).Value = 22.97892586175813
$ws.Cells.Item(11, 2).Value = 9.934245163249114
$ws.Cells.Item(11, 3).Value = 6.219938267475892
$ws.Cells.Item(11, 4).Value = 5.402097980419982
$ws.Cells.Item(11, 5).Value = 13.86081283664133
$ws.Cells.Item(11, 6).Value = 25.94895079484213
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 11).Value = 9.005148249071246
$ws.Cells.Item(11, 13).Value = 14.22194387869303
$ws.Cells.Item(11, 15).Value = 22.95883672461181
$ws.Cells.Item(12, 2).Value = 10.01309129694614
$ws.Cells.Item(12, 3).Value = 6.247418621791349
$ws.Cells.Item(12, 4).Value = 5.414067340635403
$ws.Cells.Item(12, 5).Value = 13.9303542369837
$ws.Cells.Item(12, 6).Value = 25.96107758436981
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 11).Value = 9.077603665385272
$ws.Cells.Item(12, 13).Value = 14.2781129025758
$ws.Cells.Item(12, 15).Value = 22.95203116339203
$ws.Cells.Item(13, 2).Value = 9.996155650796108
$ws.Cells.Item(13, 3).Value = 6.241514616605431
$ws.Cells.Item(13, 4).Value = 5.411494597515066
$ws.Cells.Item(13, 5).Value = 13.91537833504099
$ws.Cells.Item(13, 6).Value = 25.95841940169098
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 11).Value = 9.062058029598115
$ws.Cells.Item(13, 13).Value = 14.2660143010042
$ws.Cells.Item(13, 15).Value = 22.95346117032851
$ws.Cells.Item(14, 2).Value = 9.940752191662829
$ws.Cells.Item(14, 3).Value = 6.222205552482025
$ws.Cells.Item(14, 4).Value = 5.40308498706329
$ws.Cells.Item(14, 5).Value = 13.8665344328488
$ws.Cells.Item(14, 6).Value = 25.94992752843198
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 11).Value = 9.011135716467042
$ws.Cells.Item(14, 13).Value = 14.22656414310509
$ws.Cells.Item(14, 15).Value = 22.95826074623444
$ws.Cells.Item(15, 2).Value = 9.906684496817613
$ws.Cells.Item(15, 3).Value = 6.210336325532352
$ws.Cells.Item(15, 4).Value = 5.397919071126801
$ws.Cells.Item(15, 5).Value = 13.83661418771472
$ws.Cells.Item(15, 6).Value = 25.94486215114586
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 11).Value = 8.979772236542786
$ws.Cells.Item(15, 13).Value = 14.20240534055911
$ws.Cells.Item(15, 15).Value = 22.96130510202122
$ws.Cells.Item(16, 2).Value = 9.709491410553369
$ws.Cells.Item(16, 3).Value = 6.141687657596211
$ws.Cells.Item(16, 4).Value = 5.368091622568288
$ws.Cells.Item(16, 5).Value = 13.66517802889422
$ws.Cells.Item(16, 6).Value = 25.91792510725516
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 11).Value = 8.797424843315085
$ws.Cells.Item(16, 13).Value = 14.06408520191274
$ws.Cells.Item(16, 15).Value = 22.98035080792571
$ws.Cells.Item(17, 2).Value = 9.586823531943283
$ws.Cells.Item(17, 3).Value = 6.099022932730827
$ws.Cells.Item(17, 4).Value = 5.349597116111422
$ws.Cells.Item(17, 5).Value = 13.5600895144196
$ws.Cells.Item(17, 6).Value = 25.90326948173604
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 11).Value = 8.68325127934512
$ws.Cells.Item(17, 13).Value = 13.97938382242232
$ws.Cells.Item(17, 15).Value = 22.99346035155035
$ws.Cells.Item(18, 2).Value = 9.515668429300131
$ws.Cells.Item(18, 3).Value = 6.074285690959202
$ws.Cells.Item(18, 4).Value = 5.338888715160186
$ws.Cells.Item(18, 5).Value = 13.49969961181544
$ws.Cells.Item(18, 6).Value = 25.89552834469657
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 11).Value = 8.616743432059135
$ws.Cells.Item(18, 13).Value = 13.93073931780899
$ws.Cells.Item(18, 15).Value = 23.00152375059897
$ws.Cells.Item(19, 2).Value = 9.491475855453208
$ws.Cells.Item(19, 3).Value = 6.065876616345777
$ws.Cells.Item(19, 4).Value = 5.335250990556798
$ws.Cells.Item(19, 5).Value = 13.47926425408714
$ws.Cells.Item(19, 6).Value = 25.89302567945685
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 11).Value = 8.594082095036832
$ws.Cells.Item(19, 13).Value = 13.91428337447596
$ws.Cells.Item(19, 15).Value = 23.00434365105831
$ws.Cells.Item(20, 2).Value = 9.599944361180043
$ws.Cells.Item(20, 3).Value = 6.103585230038512
$ws.Cells.Item(20, 4).Value = 5.351573250129425
$ws.Cells.Item(20, 5).Value = 13.57127130132177
$ws.Cells.Item(20, 6).Value = 25.90475838130907
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 11).Value = 8.695492248921731
$ws.Cells.Item(20, 13).Value = 13.98839320958075
$ws.Cells.Item(20, 15).Value = 22.99201065747775
$ws.Cells.Item(21, 2).Value = 9.957053047487605
$ws.Cells.Item(21, 3).Value = 6.227885843333577
$ws.Cells.Item(21, 4).Value = 5.40555818210701
$ws.Cells.Item(21, 5).Value = 13.8808816055409
$ws.Cells.Item(21, 6).Value = 25.95239343580325
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 11).Value = 9.026128748589986
$ws.Cells.Item(21, 13).Value = 14.23815054532703
$ws.Cells.Item(21, 15).Value = 22.9568292178036
$ws.Cells.Item(22, 2).Value = 10.18462139694665
$ws.Cells.Item(22, 3).Value = 6.307262494278623
$ws.Cells.Item(22, 4).Value = 5.440182122186847
$ws.Cells.Item(22, 5).Value = 14.08320718466615
$ws.Cells.Item(22, 6).Value = 25.98962303251562
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 11).Value = 9.234543444498094
$ws.Cells.Item(22, 13).Value = 14.40167592904986
$ws.Cells.Item(22, 15).Value = 22.93851023583637
$ws.Cells.Item(23, 2).Value = 10.06371841312128
$ws.Cells.Item(23, 3).Value = 6.265072680778518
$ws.Cells.Item(23, 4).Value = 5.421764219497466
$ws.Cells.Item(23, 5).Value = 13.97524823009217
$ws.Cells.Item(23, 6).Value = 25.96919681400868
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 11).Value = 9.124020240527276
$ws.Cells.Item(23, 13).Value = 14.31438959959525
$ws.Cells.Item(23, 15).Value = 22.94785904648465
$ws.Cells.Item(24, 2).Value = 9.59401439633959
$ws.Cells.Item(24, 3).Value = 6.101523264427033
$ws.Cells.Item(24, 4).Value = 5.35068007529529
$ws.Cells.Item(24, 5).Value = 13.56621592603087
$ws.Cells.Item(24, 6).Value = 25.90408311703188
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 11).Value = 8.689960807615769
$ws.Cells.Item(24, 13).Value = 13.98431990134794
$ws.Cells.Item(24, 15).Value = 22.99266442414336
$ws.Cells.Item(25, 2).Value = 9.064799265140069
$ws.Cells.Item(25, 3).Value = 5.917621717070555
$ws.Cells.Item(25, 4).Value = 5.271312511361546
$ws.Cells.Item(25, 5).Value = 13.12740206836689
$ws.Cells.Item(25, 6).Value = 25.86017267072859
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 11).Value = 8.18988142953916
$ws.Cells.Item(25, 13).Value = 13.6313163832448
$ws.Cells.Item(25, 15).Value = 23.06098219567858
